# Applies the per-cell value updates from the "Updated cryptos list" GitHub
# Actions commit (Sat Sep 14 04:42:05 UTC 2024) to the crypto ranking sheet.
#
# A handful of "Price" cells hold numeric-looking text (e.g. "1.00", "39.81",
# "0.422") that Excel would otherwise auto-convert to a real number (dropping
# trailing zeros / becoming a Number-typed cell). Those assignments are given a
# leading apostrophe so they land back in the sheet as literal text, matching
# the original inlineStr cells. Values that are already unambiguous text (the
# "xx.xx.xx"-style big prices, the padded "  +n.nn%  " change cells, coin names
# and URLs) are assigned as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.289.66"
$ws.Range("E2").Value = "  +4.10%  "
$ws.Range("D3").Value = "2.430.35"
$ws.Range("E3").Value = "  +3.20%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'555.29"
$ws.Range("E5").Value = "  +2.24%  "
$ws.Range("D6").Value = "'139.31"
$ws.Range("E6").Value = "  +3.57%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E8").Value = "  +1.63%  "
$ws.Range("E9").Value = "  +4.46%  "
$ws.Range("D10").Value = "'5.77"
$ws.Range("E10").Value = "  +4.33%  "
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("E12").Value = "  -2.07%  "
$ws.Range("D13").Value = "'25.01"
$ws.Range("E13").Value = "  +5.41%  "
$ws.Range("D14").Value = "2.863.44"
$ws.Range("E14").Value = "  +3.12%  "
$ws.Range("D15").Value = "60.211.50"
$ws.Range("E15").Value = "  +4.07%  "
$ws.Range("E16").Value = "  +4.05%  "
$ws.Range("D17").Value = "2.431.91"
$ws.Range("E17").Value = "  +3.01%  "
$ws.Range("D18").Value = "'11.38"
$ws.Range("E18").Value = "  +6.17%  "
$ws.Range("D19").Value = "'4.42"
$ws.Range("E19").Value = "  +3.14%  "
$ws.Range("D20").Value = "'333.94"
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'65.22"
$ws.Range("E23").Value = "  +4.34%  "
$ws.Range("E24").Value = "  +3.61%  "
$ws.Range("E25").Value = "  +2.74%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").Value = "0.0₃0789"
$ws.Range("E28").Value = "  +7.14%  "
$ws.Range("E29").Value = "  +1.60%  "
$ws.Range("E30").Value = "  +3.33%  "
$ws.Range("D31").Value = "'169.39"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("E32").Value = "  +2.93%  "
$ws.Range("D33").Value = "'18.77"
$ws.Range("E33").Value = "  +2.02%  "
$ws.Range("E35").Value = "  +6.33%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").Value = "'325.32"
$ws.Range("E39").Value = "  +12.53%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").Value = "'0.422"
$ws.Range("E40").Value = "  +11.38%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "'39.81"
$ws.Range("E41").Value = "  +1.75%  "
$ws.Range("D42").Value = "'3.71"
$ws.Range("D43").Value = "'140.78"
$ws.Range("E43").Value = "  -1.38%  "
$ws.Range("E44").Value = "  +3.76%  "
$ws.Range("E45").Value = "  +1.28%  "
$ws.Range("D46").Value = "'19.65"
$ws.Range("E46").Value = "  +2.48%  "
$ws.Range("E47").Value = "  +8.88%  "
$ws.Range("E48").Value = "  +1.26%  "
$ws.Range("D49").Value = "'0.0227"
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("D50").Value = "'17.93"
$ws.Range("E50").Value = "  +2.57%  "
